$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 2
    4  = 1
    5  = 1
    6  = 8
    7  = 2
    8  = 3
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 4
    16 = 7
    17 = 3
    18 = 4
    19 = 2
    20 = 4
    21 = 1
    22 = 5
    23 = 1
    24 = 5
    25 = 2
    26 = 6
    27 = 2
    28 = 6
    29 = 0
    30 = 1
    31 = 0
    32 = 5
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
